$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '43.966.50'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '2.315.05'
$ws.Range('E3').Value = '  +1.61%  '
Set-TextValue $ws.Range('D4') '1.01'
$ws.Range('E4').Value = '  +0.34%  '
Set-TextValue $ws.Range('D5') '116.54'
$ws.Range('E5').Value = '  +23.09%  '
Set-TextValue $ws.Range('D6') '270.28'
$ws.Range('E6').Value = '  +1.34%  '
Set-TextValue $ws.Range('D7') '0.628'
$ws.Range('E7').Value = '  +1.18%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  +4.01%  '
Set-TextValue $ws.Range('D10') '49.20'
$ws.Range('E10').Value = '  +10.78%  '
Set-TextValue $ws.Range('D11') '0.0948'
$ws.Range('E11').Value = '  +1.60%  '
Set-TextValue $ws.Range('D12') '8.75'
$ws.Range('E12').Value = '  +13.63%  '
$ws.Range('E13').Value = '  +2.41%  '
Set-TextValue $ws.Range('D14') '15.75'
$ws.Range('E14').Value = '  +4.22%  '
$ws.Range('D15').Value = '2.629.75'
$ws.Range('E15').Value = '  +0.37%  '
Set-TextValue $ws.Range('D16') '0.866'
$ws.Range('E16').Value = '  +2.66%  '
$ws.Range('D17').Value = '2.323.62'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('D18').Value = '43.889.43'
$ws.Range('E18').Value = '  +1.02%  '
Set-TextValue $ws.Range('D19') '0.0000110'
$ws.Range('E19').Value = '  +3.75%  '
Set-TextValue $ws.Range('D20') '6.63'
$ws.Range('E20').Value = '  +7.43%  '
Set-TextValue $ws.Range('D21') '72.71'
$ws.Range('E21').Value = '  +1.05%  '
Set-TextValue $ws.Range('D22') '2.56'
$ws.Range('E22').Value = '  +8.38%  '
Set-TextValue $ws.Range('D23') '235.09'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D24') '9.55'
$ws.Range('E24').Value = '  +6.59%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D25') '2.90'
$ws.Range('E25').Value = '  +16.42%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +3.36%  '
Set-TextValue $ws.Range('D28') '43.88'
$ws.Range('E28').Value = '  +12.75%  '
$ws.Range('E29').Value = '  -1.29%  '
Set-TextValue $ws.Range('D30') '2.27'
$ws.Range('E30').Value = '  -0.03%  '
Set-TextValue $ws.Range('D31') '178.03'
$ws.Range('E31').Value = '  +1.60%  '
Set-TextValue $ws.Range('D32') '21.89'
$ws.Range('E32').Value = '  +0.46%  '
Set-TextValue $ws.Range('D33') '0.0938'
$ws.Range('E33').Value = '  +6.33%  '
Set-TextValue $ws.Range('D34') '5.59'
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('E35').Value = '  +8.74%  '
$ws.Range('E36').Value = '  +1.23%  '
Set-TextValue $ws.Range('D37') '0.110'
$ws.Range('E37').Value = '  +3.56%  '
Set-TextValue $ws.Range('D38') '3.97'
$ws.Range('E38').Value = '  +20.66%  '
Set-TextValue $ws.Range('D39') '0.0358'
$ws.Range('E39').Value = '  +1.59%  '
Set-TextValue $ws.Range('D40') '0.251'
$ws.Range('E40').Value = '  +6.85%  '
Set-TextValue $ws.Range('D41') '75.31'
$ws.Range('E41').Value = '  +20.36%  '
Set-TextValue $ws.Range('D42') '2.40'
$ws.Range('E42').Value = '  +2.91%  '
Set-TextValue $ws.Range('D43') '13.31'
$ws.Range('E43').Value = '  +12.28%  '
$ws.Range('E44').Value = '  +17.00%  '
$ws.Range('E45').Value = '  +0.27%  '
Set-TextValue $ws.Range('D46') '1.39'
$ws.Range('E46').Value = '  +4.65%  '
Set-TextValue $ws.Range('D47') '8.82'
$ws.Range('E47').Value = '  +0.50%  '
Set-TextValue $ws.Range('D48') '0.102'
$ws.Range('E48').Value = '  -0.05%  '
Set-TextValue $ws.Range('D49') '100.89'
$ws.Range('E49').Value = '  +3.60%  '
Set-TextValue $ws.Range('D50') '1.23'
$ws.Range('E50').Value = '  +4.33%  '
$ws.Range('E51').Value = '  +9.20%  '
